# This script applies spell-check style run-splitting (w:proofErr
# spellStart/spellEnd wrapping around words Word's spell checker flagged)
# to several paragraphs, and appends a new "talk about GIMP " paragraph
# at the end of the document (ahead of the trailing bookmark), matching
# the target diff.
#
# Because the COM shim's Range.InsertXML only reliably replaces content
# when invoked on a Range that covers one or more whole <w:p> paragraphs,
# each edit below locates the paragraph, then rewrites its full contents
# (and, for the final paragraph, appends the extra paragraphs after it)
# via InsertXML.

$d = $word.ActiveDocument
$w = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Get-ParagraphByText($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) " (friend of sarahs)" -> " (friend of " / spellStart / "sarahs" /
#    spellEnd / ")" -- the "Hey! ... Tara (friend of sarahs)" paragraph.
# ---------------------------------------------------------------------
$p1 = Get-ParagraphByText $d "friend of sarahs"
$rPr1 = '<w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:color w:val="444950"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:shd w:val="clear" w:color="auto" w:fill="F1F0F0"/></w:rPr>'
$xml1 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '>' + $rPr1 + '<w:t>Hey! We just use the curriculum online website for all our lesson plans and all the strands and strand units are on that. Last year on placement I used the Busy at Maths book there' + [char]0x2019 + 's loads of different ones you can get but I think a lot of schools use that one- Tara</w:t></w:r>' +
        '<w:r ' + $w + '>' + $rPr1 + '<w:t xml:space="preserve"> (friend of </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '>' + $rPr1 + '<w:t>sarahs</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '>' + $rPr1 + '<w:t>)</w:t></w:r>' +
        '</w:p>'
[void]$p1.Range.InsertXML($xml1)

# ---------------------------------------------------------------------
# 2) "...check the email and username in the db multiple times..."
# ---------------------------------------------------------------------
$p2 = Get-ParagraphByText $d "Talk about speed when doing verification"
$xml2 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:t xml:space="preserve">Talk about speed when doing verification by doing the password first to avoid having to check the email and username in the </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>db</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> multiple times if they get the password wrong</w:t></w:r>' +
        '</w:p>'
[void]$p2.Range.InsertXML($xml2)

# ---------------------------------------------------------------------
# 3) "To avoid rechecking the db many times..."
# ---------------------------------------------------------------------
$p3 = Get-ParagraphByText $d "To avoid rechecking the db"
$xml3 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:t xml:space="preserve">To avoid rechecking the </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>db</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> many times, when finally able to keep info between pages, keep everything loaded. E.g. for classroom, should know if the user has a classroom or not.</w:t></w:r>' +
        '</w:p>'
[void]$p3.Range.InsertXML($xml3)

# ---------------------------------------------------------------------
# 4) "Give everyone in the db an ID for reference."
# ---------------------------------------------------------------------
$p4 = Get-ParagraphByText $d "Give everyone in the db"
$xml4 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:t xml:space="preserve">Give everyone in the </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>db</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> an ID for reference.</w:t></w:r>' +
        '</w:p>'
[void]$p4.Range.InsertXML($xml4)

# ---------------------------------------------------------------------
# 5) "Overall XP only determines ... cant earn anymore overall xp form..."
#    (keeps the w:lastRenderedPageBreak run marker)
# ---------------------------------------------------------------------
$p5 = Get-ParagraphByText $d "Overall XP only determines"
$xml5 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:lastRenderedPageBreak/><w:t xml:space="preserve">Overall XP only determines what level you unlock. Once you have reached a threshold you </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>cant</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> earn anymore overall </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>xp</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> form the previous levels to stop boosting from the easy levels.</w:t></w:r>' +
        '</w:p>'
[void]$p5.Range.InsertXML($xml5)

# ---------------------------------------------------------------------
# 6) "Research children and color and use color that kids like the most"
# ---------------------------------------------------------------------
$p6 = Get-ParagraphByText $d "Research children and color"
$xml6 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:t xml:space="preserve">Research children and </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>color</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> and use </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>color</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> that kids like the most</w:t></w:r>' +
        '</w:p>'
[void]$p6.Range.InsertXML($xml6)

# ---------------------------------------------------------------------
# 7) "mention allow_copy false for " + "password" (second run untouched)
# ---------------------------------------------------------------------
$p7 = Get-ParagraphByText $d "mention allow_copy false for"
$xml7 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:t xml:space="preserve">mention </w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellStart"/>' +
        '<w:r ' + $w + '><w:t>allow_copy</w:t></w:r>' +
        '<w:proofErr ' + $w + ' w:type="spellEnd"/>' +
        '<w:r ' + $w + '><w:t xml:space="preserve"> false for </w:t></w:r>' +
        '<w:r ' + $w + '><w:t>password</w:t></w:r>' +
        '</w:p>'
[void]$p7.Range.InsertXML($xml7)

# ---------------------------------------------------------------------
# 8) Final paragraph: keep "...add it to the docume" + "nt" runs (minus
#    the bookmark), add a blank paragraph, then a new paragraph holding
#    "talk about GIMP " followed by the _GoBack bookmark pair.
# ---------------------------------------------------------------------
$p8 = Get-ParagraphByText $d "add confirm password"
$xml8 = '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:t>add confirm password and look up why this is a thing and add it to the docume</w:t></w:r>' +
        '<w:r ' + $w + '><w:t>nt</w:t></w:r>' +
        '</w:p>' +
        '<w:p ' + $w + '/>' +
        '<w:p ' + $w + '>' +
        '<w:r ' + $w + '><w:t xml:space="preserve">talk about GIMP </w:t></w:r>' +
        '<w:bookmarkStart ' + $w + ' w:id="0" w:name="_GoBack"/>' +
        '<w:bookmarkEnd ' + $w + ' w:id="0"/>' +
        '</w:p>'
[void]$p8.Range.InsertXML($xml8)

Write-Output "Applied all edits"
